# Auto-generated edit script applying scheduled-runner value updates
# to the Chocobo_Profits workbook's per-job sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1354.8507
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 4962.5
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 14887.5
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -15427.5

$ws.Range("H73").Value = 1354.8507
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 4962.5
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 14887.5
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -16759.5

$ws.Range("H115").Value = 1220.9
$ws.Range("I115").Value = 1220.9
$ws.Range("K115").Value = 3662.7
$ws.Range("M115").Value = -2095.7

$ws.Range("H133").Value = 37785.715
$ws.Range("J133").Value = 37785.715
$ws.Range("L133").Value = 37785.715
$ws.Range("N133").Value = -47905.715

$ws.Range("H141").Value = 3715.8333
$ws.Range("I141").Value = 3473.75
$ws.Range("K141").Value = 10421.25
$ws.Range("M141").Value = -5241.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2760.7144
$ws.Range("I74").Value = 2112.375
$ws.Range("K74").Value = 2112.375
$ws.Range("M74").Value = -1238.375

$ws.Range("H77").Value = 2760.7144
$ws.Range("I77").Value = 2112.375
$ws.Range("K77").Value = 10561.875
$ws.Range("M77").Value = -6193.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13160941
$ws.Range("I31").Value = 1820.5714
$ws.Range("J31").Value = 50006480
$ws.Range("K31").Value = 1820.5714
$ws.Range("L31").Value = 50006480
$ws.Range("M31").Value = -1525.5714
$ws.Range("N31").Value = -50007070

$ws.Range("H34").Value = 13160941
$ws.Range("I34").Value = 1820.5714
$ws.Range("J34").Value = 50006480
$ws.Range("K34").Value = 1820.5714
$ws.Range("L34").Value = 50006480
$ws.Range("M34").Value = -1618.5714
$ws.Range("N34").Value = -50006884

$ws.Range("H39").Value = 15276.272
$ws.Range("I39").Value = 3996
$ws.Range("J39").Value = 24676.5
$ws.Range("K39").Value = 3996
$ws.Range("L39").Value = 24676.5
$ws.Range("M39").Value = -3605
$ws.Range("N39").Value = -25458.5

$ws.Range("H44").Value = 23866.334
$ws.Range("J44").Value = 23866.334
$ws.Range("L44").Value = 23866.334
$ws.Range("N44").Value = -24750.334

$ws.Range("H49").Value = 15276.272
$ws.Range("I49").Value = 3996
$ws.Range("J49").Value = 24676.5
$ws.Range("K49").Value = 3996
$ws.Range("L49").Value = 24676.5
$ws.Range("M49").Value = -3814
$ws.Range("N49").Value = -25040.5

$ws.Range("H57").Value = 49731.57
$ws.Range("J57").Value = 49731.57
$ws.Range("L57").Value = 49731.57
$ws.Range("N57").Value = -50851.57

$ws.Range("H107").Value = 647.2381
$ws.Range("I107").Value = 410.70587
$ws.Range("K107").Value = 410.70587
$ws.Range("M107").Value = 1509.29413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 290.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 290.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 871.5
$ws.Range("N23").Value = -1341.5
$ws.Range("M23").ClearContents()

$ws.Range("H41").Value = 1743.3334
$ws.Range("J41").Value = 2355
$ws.Range("L41").Value = 7065
$ws.Range("N41").Value = -7741

$ws.Range("H131").Value = 6945450.5
$ws.Range("I131").Value = 166670140
$ws.Range("J131").Value = 898.1449
$ws.Range("K131").Value = 500010420
$ws.Range("L131").Value = 2694.4347
$ws.Range("M131").Value = -500005380
$ws.Range("N131").Value = -12774.4347

$ws.Range("H132").Value = 1994.5
$ws.Range("I132").Value = 866.6667
$ws.Range("J132").Value = 2671.2
$ws.Range("K132").Value = 7800.0003
$ws.Range("L132").Value = 24040.8
$ws.Range("M132").Value = -5270.0003
$ws.Range("N132").Value = -29100.8

$ws.Range("H137").Value = 2759.05
$ws.Range("J137").Value = 2980.611
$ws.Range("L137").Value = 8941.832999999999
$ws.Range("N137").Value = -19141.833

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2661.3667
$ws.Range("I122").Value = 1563.8334
$ws.Range("K122").Value = 4691.5002
$ws.Range("M122").Value = -2241.5002

$ws.Range("H123").Value = 10794.071
$ws.Range("J123").Value = 10794.071
$ws.Range("L123").Value = 10794.071
$ws.Range("N123").Value = -15694.071

$ws.Range("H132").Value = 6479.8
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2320.8235
$ws.Range("I22").Value = 1164.5714
$ws.Range("K22").Value = 1164.5714
$ws.Range("M22").Value = -869.5714

$ws.Range("H27").Value = 2320.8235
$ws.Range("I27").Value = 1164.5714
$ws.Range("K27").Value = 1164.5714
$ws.Range("M27").Value = -1057.5714

$ws.Range("H45").Value = 28221.666
$ws.Range("I45").Value = 29000
$ws.Range("J45").Value = 27832.5
$ws.Range("K45").Value = 29000
$ws.Range("L45").Value = 27832.5
$ws.Range("N45").Value = -28646.5
$ws.Range("M45").Value = -28593

$ws.Range("H46").Value = 1777.1482
$ws.Range("I46").Value = 1580.0667
$ws.Range("J46").Value = 2023.5
$ws.Range("K46").Value = 1580.0667
$ws.Range("L46").Value = 2023.5
$ws.Range("M46").Value = -1392.0667
$ws.Range("N46").Value = -2399.5

$ws.Range("H47").Value = 39999
$ws.Range("J47").Value = 39999
$ws.Range("L47").Value = 39999
$ws.Range("N47").Value = -40979

$ws.Range("H52").Value = 39999
$ws.Range("J52").Value = 39999
$ws.Range("L52").Value = 39999
$ws.Range("N52").Value = -40465

$ws.Range("H68").Value = 760.26666
$ws.Range("I68").Value = 709.36365
$ws.Range("K68").Value = 709.36365
$ws.Range("M68").Value = 39.63634999999999

$ws.Range("H71").Value = 760.26666
$ws.Range("I71").Value = 709.36365
$ws.Range("K71").Value = 3546.81825
$ws.Range("M71").Value = 197.1817499999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100022520
$ws.Range("I62").Value = 166669260
$ws.Range("J62").Value = 52400
$ws.Range("K62").Value = 166669260
$ws.Range("L62").Value = 52400
$ws.Range("M62").Value = -166668636
$ws.Range("N62").Value = -53648

$ws.Range("H65").Value = 100022520
$ws.Range("I65").Value = 166669260
$ws.Range("J65").Value = 52400
$ws.Range("K65").Value = 833346300
$ws.Range("L65").Value = 262000
$ws.Range("M65").Value = -833343180
$ws.Range("N65").Value = -268240
